# Daily attendance processing - 2026-01-07 05:15:56
#
# Column G ("Recorded By") holds a comma-separated list of the people/
# processes that touched the attendance record. For every row whose list
# has more than one entry and does NOT already start with "System", the
# order of the entries is reversed so that "System" (when present) moves
# to the front, e.g. "backup@backdoor.com, System" -> "System, backup@backdoor.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 157
$col = 7  # column G

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $v = $cell.Value2

    if ($v -eq $null) { continue }
    if ($v -eq "") { continue }

    $parts = $v -split ", "

    # Use .Equals() (not -eq/-ne) because PowerShell's default comparison
    # operators are case-insensitive, and "system" must be treated as
    # different from "System" here.
    $startsWithSystem = $parts[0].Equals("System")

    if ($parts.Length -gt 1 -and -not $startsWithSystem) {
        $n = $parts.Length
        $reversed = @()
        for ($i = $n - 1; $i -ge 0; $i--) {
            $reversed += $parts[$i]
        }
        $newVal = $reversed -join ", "
        $cell.Value = $newVal
    }
}
